$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Ngày" (Date) column header in the revenue export header row to
# "Thời gian" (Time), per "Fix branch revenue" task.
$ws.Range("B8").Value = "Thời gian"

# Update the last active/selected cell as left by the author when saving.
$ws.Range("E18").Select()
